$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update F2 and F3 from "v0.99.1dev" to "NewDriver5 2023.12.1.mgu"
$ws.Range("F2").Value = "NewDriver5 2023.12.1.mgu"
$ws.Range("F3").Value = "NewDriver5 2023.12.1.mgu"

# Insert new note text into row 16, column G
$ws.Range("G16").Value = "Note that v0.99.1dev has poor support for high-rate PSG music (Atart ST SNDH VGM dumps)."

# Apply same style as the other remark cells (row height + wrap)
$ws.Rows.Item(16).RowHeight = 42.75
$ws.Range("G16").WrapText = $true
$ws.Range("G16").HorizontalAlignment = -4131
$ws.Range("G16").VerticalAlignment = -4108

# Update selection to G16
$ws.Range("G16").Select()
